$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 156.9520094073222
$ws.Cells.Item(3, 3).Value = 11.54891641892396
$ws.Cells.Item(4, 3).Value = 9.028272274071071
$ws.Cells.Item(5, 3).Value = 14.20144173836713
$ws.Cells.Item(6, 3).Value = 34.17510641205065
$ws.Cells.Item(7, 3).Value = 10.6719438624676
$ws.Cells.Item(8, 3).Value = 7.257935150757307
$ws.Cells.Item(9, 3).Value = 26.32929416100103
$ws.Cells.Item(10, 3).Value = 47.2916950984208
$ws.Cells.Item(11, 3).Value = 8.958978775557526
$ws.Cells.Item(12, 3).Value = 3.469890565350249
$ws.Cells.Item(13, 3).Value = 6.818331235456331
$ws.Cells.Item(14, 3).Value = 1.566927176064328
$ws.Cells.Item(15, 3).Value = 1.360536863287429
$ws.Cells.Item(16, 3).Value = 21.04808644633438
$ws.Cells.Item(17, 3).Value = 21.7961581937279
$ws.Cells.Item(18, 3).Value = 9.661599948657223
$ws.Cells.Item(19, 3).Value = 1.080382503705789
$ws.Cells.Item(20, 3).Value = 22.72454205553302
$ws.Cells.Item(21, 3).Value = 66.50387637983721
$ws.Cells.Item(22, 3).Value = 7.466560737679805
$ws.Cells.Item(23, 3).Value = 2.393978609935656
$ws.Cells.Item(24, 3).Value = 24.11413748271322
$ws.Cells.Item(25, 3).Value = 5.363167766671912
$ws.Cells.Item(26, 3).Value = 9.620619922654589
$ws.Cells.Item(27, 3).Value = 19.91852791142543
$ws.Cells.Item(28, 3).Value = 7.309346456106066
$ws.Cells.Item(29, 3).Value = 5.243208054191476
$ws.Cells.Item(30, 3).Value = 3.098089965799085
$ws.Cells.Item(31, 3).Value = 1.756180387058308
$ws.Cells.Item(32, 3).Value = 4.753683016305474
$ws.Cells.Item(33, 3).Value = 2.233783962834453
$ws.Cells.Item(34, 3).Value = 88.93112697400565
$ws.Cells.Item(35, 3).Value = 8.171417184925099
$ws.Cells.Item(36, 3).Value = 24.04633416696341
$ws.Cells.Item(37, 3).Value = 4.115884793428125
$ws.Cells.Item(38, 3).Value = 3.680006335036479
$ws.Cells.Item(39, 3).Value = 8.631138567536459
$ws.Cells.Item(40, 3).Value = 0.9753246188626747
$ws.Cells.Item(41, 3).Value = 5.531558418973643
$ws.Cells.Item(42, 3).Value = 280.2
